$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 194
$ws.Range("F5").Value = 1735
$ws.Range("F9").Value = 3554
$ws.Range("F10").Value = 932
$ws.Range("F11").Value = 1178
$ws.Range("F12").Value = 1580
$ws.Range("G12").Value = '已售罄'
$ws.Range("F13").Value = 29
$ws.Range("F14").Value = 888
$ws.Range("F16").Value = 1331
$ws.Range("F17").Value = 1793
$ws.Range("F19").Value = 461
$ws.Range("F20").Value = 1550
$ws.Range("F21").Value = 13
$ws.Range("F22").Value = 1079
$ws.Range("F23").Value = 2314
$ws.Range("F24").Value = 16
$ws.Range("F26").Value = 4314
$ws.Range("F28").Value = 6
$ws.Range("F31").Value = 1218

$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 18
$ws.Range("F10").Value = 165
$ws.Range("F22").Value = 133
$ws.Range("F24").Value = 192
$ws.Range("F39").Value = 26

$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 9578
$ws.Range("F9").Value = 391
$ws.Range("F10").Value = 2982
$ws.Range("F11").Value = 501
$ws.Range("F12").Value = 810
$ws.Range("F13").Value = 219
$ws.Range("F14").Value = 247

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 194
$ws.Range("B6").Value = '2024-09-24'
$ws.Range("C6").Value = '上海·星零界·社交游乐·休闲运动·潮玩派对'
$ws.Range("D6").Value = '长宁路1191号长宁来福士B1 长宁来福士'
$ws.Range("E6").Value = '2024.09.24 10:00-12.31 22:00'
$ws.Range("F6").Value = 16
$ws.Range("G6").Value = 68
$ws.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=92659'
$ws.Range("I6").Value = '//i0.hdslb.com/bfs/openplatform/202409/PHS8s1lu1726221065737.png'
$ws.Range("B7").Value = '2024-09-26'
$ws.Range("C7").Value = '上海·【神秘的西夏陵】大空间高沉浸探险体验'
$ws.Range("D7").Value = '南京西路325号 上海市历史博物馆'
$ws.Range("E7").Value = '2024.09.26 10:00-12.31 19:00'
$ws.Range("F7").Value = 33
$ws.Range("G7").Value = 108
$ws.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=92581'
$ws.Range("I7").Value = '//i2.hdslb.com/bfs/openplatform/202409/jeDZO2cS1726302714881.jpeg'
$ws.Range("B8").Value = '2024-09-28'
$ws.Range("C8").Value = '上海·［咒术回战 2024 剧场版 咒术回战 0］主题咖啡厅'
$ws.Range("D8").Value = '大悦城 次元波板糖'
$ws.Range("E8").Value = '2024.09.28 00:00-10.27 23:59'
$ws.Range("F8").Value = 391
$ws.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=92608'
$ws.Range("I8").Value = '//i2.hdslb.com/bfs/openplatform/202409/DBTiL9sY1726727259104.png'
$ws.Range("C9").Value = '上海·2024·《世界之外》x  萌果酱谷子咖啡'
$ws.Range("D9").Value = '南京东路340号百联ZX 萌果酱谷子咖啡（百联）'
$ws.Range("E9").Value = '2024.10.01 00:00-12.11 23:59'
$ws.Range("F9").Value = 2982
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=93006'
$ws.Range("I9").Value = '//i2.hdslb.com/bfs/openplatform/202409/qtffZOKB1727426243733.png'
$ws.Range("B10").Value = '2024-10-01'
$ws.Range("C10").Value = '上海·三丽鸥家族Sanrio Characters主题餐厅·海滩奇遇季'
$ws.Range("D10").Value = '南京东路800号4楼 上海市第一百货商店-C馆'
$ws.Range("E10").Value = '2024.10.01 00:00-11.19 23:59'
$ws.Range("F10").Value = 501
$ws.Range("G10").Value = 10
$ws.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=93078'
$ws.Range("I10").Value = '//i0.hdslb.com/bfs/openplatform/202409/aiu4g5K21727677592777.png'
$ws.Range("B11").Value = '2024-10-10'
$ws.Range("C11").Value = '上海·「火影忍者疾风传 × animate cafe」'
$ws.Range("D11").Value = '西藏北路198号大悦城北座8楼N809-1 animate cafe上海店'
$ws.Range("E11").Value = '2024.10.10 00:00-11.12 23:59'
$ws.Range("F11").Value = 810
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=92883'
$ws.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202409/aQIhaIgt1727249498713.png'
$ws.Range("B12").Value = '2024-10-17'
$ws.Range("C12").Value = '上海·蜡笔小新：我们的恐龙日记x HAPPY ZOO 主题咖啡厅'
$ws.Range("D12").Value = '南京东路340号百联zx创趣场四楼05号 HAPPY ZOO'
$ws.Range("E12").Value = '2024.10.17 00:00-10.27 23:59'
$ws.Range("F12").Value = 219
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=93221'
$ws.Range("I12").Value = '//i0.hdslb.com/bfs/openplatform/202410/nzGP5KRA1728526131597.png'
$ws.Range("C13").Value = '上海·ROOKiEZ is PUNK`D 「Reignite Youth （重燃青春）」2024 CHINA Tour '
$ws.Range("D13").Value = '虹许路731号4号楼 THE BOXX•城市乐园'
$ws.Range("E13").Value = '2024.10.18 20:30-10.18 22:00'
$ws.Range("F13").Value = 81
$ws.Range("G13").Value = 259
$ws.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=91376'
$ws.Range("I13").Value = '//i1.hdslb.com/bfs/openplatform/202408/pZdI02BJ1724735899119.jpeg'
$ws.Range("C14").Value = '上海·“爆裂鼓手”电影中的鼓手经典音乐会'
$ws.Range("D14").Value = '南京西路1376号 上海商城剧院'
$ws.Range("E14").Value = '2024.10.18 19:30-10.18 21:00'
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 100
$ws.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=93194'
$ws.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202410/b0buY3ha1728359625532.jpeg'
$ws.Range("B15").Value = '2024-10-18'
$ws.Range("C15").Value = '上海·幻境尘影—无期迷途光影展'
$ws.Range("D15").Value = '陆家嘴西路168号 上海正大广场'
$ws.Range("E15").Value = '2024.10.18 00:00-11.30 23:59'
$ws.Range("F15").Value = 1735
$ws.Range("G15").Value = 98
$ws.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=93269'
$ws.Range("I15").Value = '//i1.hdslb.com/bfs/openplatform/202410/kq6Pmpkv1728728206204.png'
$ws.Range("C16").Value = '上海·排球少年Only·魔都见学同人展'
$ws.Range("D16").Value = '吴中路1588号上海爱琴海购物中心F4 竞梦元宇宙'
$ws.Range("E16").Value = '2024.10.19 10:00-10.20 18:00'
$ws.Range("F16").Value = 320
$ws.Range("G16").Value = 89
$ws.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=92129'
$ws.Range("I16").Value = '//i2.hdslb.com/bfs/openplatform/202409/NpLGFHMN1725614230335.jpeg'
$ws.Range("B17").Value = '2024-10-19'
$ws.Range("C17").Value = '上海·第五人格同人only'
$ws.Range("D17").Value = '漕宝路1688号 诺宝中心酒店'
$ws.Range("E17").Value = '2024.10.19 10:00-10.19 17:00'
$ws.Range("F17").Value = 619
$ws.Range("G17").Value = 60
$ws.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=91590'
$ws.Range("I17").Value = '//i2.hdslb.com/bfs/openplatform/202409/XE1wxOQn1725446984161.jpeg'
$ws.Range("C18").Value = '上海·LookLook动漫嘉年华3th'
$ws.Range("D18").Value = '曹安公路4218号 上海国际短视频中心'
$ws.Range("E18").Value = '2024.10.26 10:00-10.27 17:30'
$ws.Range("F18").Value = 932
$ws.Range("G18").Value = 68
$ws.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=90495'
$ws.Range("I18").Value = '//i0.hdslb.com/bfs/openplatform/202409/DQLGW65C1726814328151.jpeg'
$ws.Range("C19").Value = '上海·【早鸟4折】“海上钢琴师”一生必听经典电影主题音乐会'
$ws.Range("D19").Value = '南京西路1376号 上海商城剧院'
$ws.Range("E19").Value = '2024.10.26 15:00-10.26 16:30'
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=91375'
$ws.Range("I19").Value = '//i1.hdslb.com/bfs/openplatform/202408/qfPgppOK1724743485013.jpeg'
$ws.Range("F22").Value = 165
$ws.Range("F27").Value = 1793
$ws.Range("F29").Value = 1550
$ws.Range("F30").Value = 133
$ws.Range("F31").Value = 133
$ws.Range("F32").Value = 13
$ws.Range("F33").Value = 192
$ws.Range("F34").Value = 1079
$ws.Range("F36").Value = 2314
$ws.Range("F40").Value = 247
$ws.Range("F42").Value = 6
$ws.Range("F44").Value = 26
$ws.Range("F50").Value = 1218
